$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("System Overview")
$wsOverview.Activate()
$wsOverview.Range("D7:M7").Select()

$wsHardware = $wb.Worksheets.Item("Hardware Lifecycles")
$wsHardware.Activate()
$wsHardware.Range("A2:G2").Select()

$wsSystemData = $wb.Worksheets.Item("System Data")
$wsSystemData.Activate()
$wsSystemData.Range("D34").Select()

$before = $wb.Worksheets.Item("System Data")
$wsNew = $wb.Worksheets.Add($before)
$wsNew.Name = "Modernization Timeline"
$wsNew.Activate()
